$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.132.24"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "3.444.39"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'412.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").Value = "'130.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "

# Row 7
$ws.Range("D7").Value = "'0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.64%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  -2.51%  "

# Row 10
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("E12").Value = "  +14.21%  "

# Row 13
$ws.Range("D13").Value = "'9.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.67%  "

# Row 14
$ws.Range("D14").Value = "3.983.02"
$ws.Range("E14").Value = "  -0.04%  "

# Row 16
$ws.Range("D16").Value = "'21.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("D17").Value = "3.430.53"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("D18").Value = "'12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("E19").Value = "  +2.30%  "

# Row 20
$ws.Range("D20").Value = "62.127.82"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").Value = "'502.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +23.73%  "

# Row 22
$ws.Range("D22").Value = "'93.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.64%  "

# Row 23
$ws.Range("E23").Value = "  +3.60%  "

# Row 24
$ws.Range("E24").Value = "  +1.82%  "

# Row 25
$ws.Range("D25").Value = "'3.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.82%  "

# Row 26
$ws.Range("D26").Value = "'35.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.96%  "

# Row 27
$ws.Range("E27").Value = "  +6.55%  "

# Row 28
$ws.Range("E28").Value = "  +0.33%  "

# Row 29
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").Value = "'12.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.21%  "

# Row 31
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("E32").Value = "  -2.04%  "

# Row 33
$ws.Range("E33").Value = "  -2.19%  "

# Row 34
$ws.Range("E34").Value = "  -4.10%  "

# Row 35
$ws.Range("D35").Value = "'60.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.10%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").Value = "'0.0501"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "

# Row 38
$ws.Range("E38").Value = "  +2.77%  "

# Row 39
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("E40").Value = "  +5.06%  "

# Row 41
$ws.Range("D41").Value = "'2.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.22%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'2.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.96%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'148.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.42%  "

# Row 44
$ws.Range("D44").Value = "'2.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "

# Row 45
$ws.Range("E45").Value = "  +1.82%  "

# Row 46
$ws.Range("E46").Value = "  +7.02%  "

# Row 47
$ws.Range("D47").Value = "'16.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "

# Row 48
$ws.Range("D48").Value = "'2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +21.54%  "

# Row 49
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'120.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +24.01%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'23.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.42%  "

# Row 51
$ws.Range("E51").Value = "  +19.86%  "
